$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("01-07-2021", 31814),
    @("02-07-2021", 33273),
    @("05-07-2021", 31438),
    @("06-07-2021", 30841),
    @("07-07-2021", 31069),
    @("08-07-2021", 31529),
    @("09-07-2021", 34707),
    @("12-07-2021", 36496),
    @("13-07-2021", 36041),
    @("14-07-2021", 34922),
    @("15-07-2021", 29340),
    @("19-07-2021", 28737),
    @("20-07-2021", 28478),
    @("21-07-2021", 27673),
    @("22-07-2021", 26816),
    @("23-07-2021", 26330),
    @("26-07-2021", 25454),
    @("27-07-2021", 24253),
    @("28-07-2021", 24428),
    @("29-07-2021", 24524),
    @("30-07-2021", 21743)
)

$startRow = 126
$endRow = $startRow + $data.Length - 1

# Force column A to text before writing, so strings that look like dates
# (e.g. "01-07-2021") are not silently converted into date serials by
# Excel's "looks like a date" auto-conversion.
$dateRange = $ws.Range("A$startRow`:A$endRow")
$dateRange.NumberFormat = "@"

$row = $startRow
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row = $row + 1
}

# Restore the default (unstyled) look so the new rows match the rest of the
# sheet instead of keeping the temporary text format.
$dateRange.Style = "Normal"
